# Edit the TMF8801 register map sheet:
#  - Remove the duplicate "APPREV_MINOR" row (row 5), shifting all subsequent
#    rows up by one.
#  - Correct the bit-width / high-bit-index columns (D/E) for 8-bit registers
#    that were incorrectly listed as 16 bits wide (D=16,E=15 -> D=8,E=7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TMF8801")
$ws.Activate()

# Remove the duplicate row (old row 5: APPREV_MINOR @ 0x12), shifting rows up.
$ws.Rows("5").Delete()

# After the deletion, the data now spans rows 2..42. Any row whose bit-width
# (column D) is still 16 with a high-bit-index (column E) of 15 should really
# be an 8-bit register with high-bit-index 7.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $width = $ws.Cells.Item($r, 4).Value()
    $high = $ws.Cells.Item($r, 5).Value()
    if ($width -eq 16 -and $high -eq 15) {
        $ws.Cells.Item($r, 4).Value = 8
        $ws.Cells.Item($r, 5).Value = 7
    }
}

# Restore the cursor/selection position recorded for this sheet.
$ws.Range("E13").Select()

